$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 446
$range = $ws.Range("C2:C$lastRow")
$range.Value = 46075
